$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  'Closures, remote learning and relocationsGoogle Tag Manager (noscript)End Google Tag Manager (noscript)You may be trying to access this site from a secured browser on the server. Please enable scripts and reload this page.Start of Site headerSkip to contentSearchSearchEnd of Site headerHome>The Department>Program directory>Emergencies and Natural Disasters>Closures, remote learning and relocationsClosures, remote learning and relocations- Example page content area using current sharepoint wrapperPage Content',
  'li: This page lists early childhood services, schools and TAFEs that are currently closed.',
  'li: Closures may be the result of a number of factors including: the Covid-19 pandemic, bushfire preparedness, construction or other factors. More information is available from your local service, school or TAFE.',
  'li: See',
  'Coronavirus advice',
  'li: for the latest advice on the outbreak of COVID-19 (coronavirus).',
  'On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school and early childhood service, TAFE closures and relocations for Tuesday 8 September, (as at 10:30am, 8 September)South-Eastern Victoria RegionEarly childhood services',
  'li: The Department has been advised of the following early childhood service closures:',
  'li: Aqua Energy Creche SALE',
  'li: Bairnsdale Aquatic & Recreation Centre BAIRNSDALE',
  'li: Barring Djinang Kindergarten SOUTHBANK',
  'li: Berwick Neighbourhood Centre - Marriott Waters LYNDHURST',
  'li: Berwick Neighbourhood Centre (Timbarra Playroom) BERWICK',
  'li: Botany Park Preschool CARRUM DOWNS',
  'li: Brighton Grammar OSHClub BRIGHTON',
  'li: Brighton Playroom Association BRIGHTON',
  'li: Camp Australia - Aspendale Primary School OSHC ASPENDALE',
  'li: Camp Australia - East Bentleigh Primary School OSHC BENTLEIGH EAST',
  'li: Camp Australia - Geelong Grammar School Toorak Campus OSHC TOORAK',
  'li: Camp Australia - Haileybury Edrington OSHC BERWICK',
  'li: Camp Australia - Haileybury Newlands OSHC KEYSBOROUGH',
  'li: Camp Australia - Keysborough Primary School OSHC SPRINGVALE SOUTH',
  'li: Camp Australia - Loreto College Mandeville Hall OSHC TOORAK',
  'li: Camp Australia - Mentone Girls'' Grammar School OSHC MENTONE',
  'li: Camp Australia - Mentone Grammar OSHC MENTONE',
  'li: Camp Australia - Peninsula Grammar OSHC MOUNT ELIZA',
  'li: Camp Australia - Sandringham Primary School OSHC SANDRINGHAM',
  'li: Camp Australia - Sandringham Primary School OSHC SANDRINGHAM',
  'li: Camp Australia - Sirius College - Keysborough OSHC KEYSBOROUGH',
  'li: Camp Australia - St Clare''s Primary School Pakenham OSHC Officer',
  'li: Camp Australia - St Kevins College Toorak OSHC TOORAK',
  'li: Camp Australia - Toorak College OSHC MOUNT ELIZA',
  'li: Camp Australia - Woodleigh Primary School OSHC FRANKSTON',
  'li: Camp Australia - Woodleigh School Penbank Campus OSHC MOOROODUC',
  'li: Cardinia Life PAKENHAM',
  'li: Carnegie Occasional Care Centre CARNEGIE',
  'li: Castlefield Community Centre HAMPTON',
  'li: Club V Kids WARRAGUL',
  'li: Cranbourne Community House CRANBOURNE',
  'li: David Collings Centre MORNINGTON',
  'li: Fernwood Fitness Clayton CLAYTON',
  'li: Fernwood Fitness Narre Warren NARRE WARREN',
  'li: GEKA Bentleigh Kindergarten BENTLEIGH',
  'li: GEKA Caulfield South Kindergarten CAULFIELD SOUTH',
  'li: GEKA Centre Road Kindergarten BENTLEIGH EAST',
  'li: GEKA Glover Street Kindergarten BENTLEIGH EAST',
  'li: GEKA McKinnon Kindergarten MCKINNON',
  'li: GEKA Murrumbeena Kindergarten MURRUMBEENA',
  'li: GEKA Orrong Road Kindergarten ELSTERNWICK',
  'li: Genesis Traralgon TRARALGON',
  'li: GESAC Childcare BENTLEIGH EAST',
  'li: Glen Iris Creche GLEN IRIS',
  'li: Goodlife Chelsea Heights CHELSEA HEIGHTS',
  'li: Goodlife Fountain Gate NARRE WARREN',
  'li: Goodlife Karingal KARINGAL',
  'li: Hallam Community Centre Inc HALLAM',
  'li: Hoa Nghiem Primary School SPRINGVALE SOUTH',
  'li: Lakes Aquadome Creche LAKES ENTRANCE',
  'li: Little Beacons Learning Centre - Berwick BERWICK',
  'li: Little Pelicans - Pelican Park Recreation Centre HASTINGS',
  'li: Merkaz Bentleigh BENTLEIGH',
  'li: Moe Neighbourhood House MOE',
  'li: Paisley Park Early Learning Centre Cranbourne CRANBOURNE',
  'li: Pavillion - Frankston & District Netball Association Inc FRANKSTON',
  'li: Peninsula Aquatic Recreation Centre FRANKSTON',
  'li: Play Zone - Highett HIGHETT',
  'li: Playtime Warragul WARRAGUL',
  'li: Sacre Coeur Oshclub GLEN IRIS',
  'li: Sandybeach Centre SANDRINGHAM',
  'li: Singleton Preschool ENDEAVOUR HILLS',
  'li: Smaland Springvale Ikea SPRINGVALE',
  'li: Somerville Recreation Centre Child Care SOMERVILLE',
  'li: Toorak Primary School OSHClub TOORAK',
  'li: Upper Beaconsfield Community Early Learning Centre BEACONSFIELD UPPER',
  'Schools closedThe Department has been advised of the following school closures:',
  'li: St Clare''s Primary School, OFFICER',
  'TAFETheDepartment hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood services',
  'li: The Department has been advised of the following early childhood service closures:',
  'li: Aqualink Box Hill Creche BOX HILL',
  'li: Aqualink Nunawading Creche FOREST HILL',
  'li: Aquarena Childcare DONCASTER',
  'li: Billanook College - Early Learning Program MOOROOLBARK',
  'li: Camp Australia - Balwyn North Primary School OSHC BALWYN NORTH',
  'li: Camp Australia - Boronia K-12 College OSHC BORONIA',
  'li: Camp Australia - Glendal Primary School OSHC GLEN WAVERLEY',
  'li: Camp Australia - Mountain Gate Primary School OSHC FERNTREE GULLY',
  'li: Camp Australia - St Dominics Primary School OSHC CAMBERWELL EAST',
  'li: Camp Australia - St Jude the Apostle School OSHC SCORESBY',
  'li: Camp Australia - St Timothy''s School Vermont OSHC VERMONT',
  'li: Camp Australia - Strathcona Baptist Girls Junior School OSHC CANTERBURY',
  'li: Camp Australia - Trinity Grammar OSHC KEW',
  'li: Camp Australia - Yarra Valley Grammar School OSHC RINGWOOD',
  'li: Carey OSHClub Kew KEW',
  'li: Cathedral College Wangaratta Outside School Hours Care WANGARATTA',
  'li: Croydon Central Kindergarten CROYDON',
  'li: Croydon Hills Primary School OSHC Program CROYDON HILLS',
  'li: Ferntree Gully 3 Year Old Preschool FERNTREE GULLY',
  'li: Fitness First Doncaster (Playzone) DONCASTER',
  'li: Flamingo Community Group WANTIRNA SOUTH',
  'li: Genesis Creche Wantirna BORONIA',
  'li: Hawthorn Early Years HAWTHORN',
  'li: Highmount Kindergarten MOUNT WAVERLEY',
  'li: JOYFUL LEARNING EARLY LEARNING CENTRE PTY LTD MITCHAM',
  'li: Kew Neighbourhood Learning Centre KEW',
  'li: Kidz Biz RINGWOOD',
  'li: Kingswood College OSHClub BOX HILL',
  'li: Knox Leisureworks BORONIA',
  'li: Koori First Steps Preschool WODONGA',
  'li: MakerDojo HAWTHORN',
  'li: Monbulk Aquatic Centre MONBULK',
  'li: Playstation CANTERBURY',
  'li: St Andrews Christian College Outside School Hours Care WANTIRNA SOUTH',
  'li: Toolamba Outside School Hours Care TOOLAMBA',
  'li: VERMONT CHILDREN''S CENTRE VERMONT',
  'li: Wodonga Sports and Leisure Centre WEST WODONGA',
  'li: Yarra Centre YARRA JUNCTION',
  'Schools closed',
  'li: The Department has',
  'not',
  'li: been advised of any school closures.',
  'TAFE',
  'li: The Department has not been advised of any TAFE closures.',
  'North-Western Victoria RegionEarly childhood services',
  'li: The Department has been advised of the following early childhood service closures:',
  'li: Annadale Community Centre MICKLEHAM',
  'li: Bank Street Preschool CRAIGIEBURN',
  'li: Bendigo Community Children''s Services Heathcote OSHC HEATHCOTE',
  'li: Bluebird Way Child and Family Centre ROXBURGH PARK',
  'li: Bright Stars Early Years Child Care Centre EPPING',
  'li: Camp Australia - Alphington Grammar School OSHC ALPHINGTON',
  'li: Camp Australia - Ivanhoe Grammar OSHC IVANHOE',
  'li: Camp Australia - Ivanhoe Grammar Plenty Valley Campus OSHC MERNDA',
  'li: Camp Australia - Thomastown Meadows Primary School OSHC THOMASTOWN',
  'li: Camp Australia - Wattle Glen Primary School OSHC WATTLE GLEN',
  'li: Castlemaine Steiner Kindergarten MUCKLEFORD',
  'li: Collingwood College Afterschool Care and Vacation Care Program COLLINGWOOD',
  'li: Dawson Street Preschool TULLAMARINE',
  'li: Diamond Creek Community Centre DIAMOND CREEK',
  'li: Diamond Valley Sports and Fitness Centre GREENSBOROUGH',
  'li: Dobell Avenue Preschool SUNBURY',
  'li: Edwardes Street Child Minding RESERVOIR',
  'li: Eltham Leisure Centre ELTHAM',
  'li: Greensborough Preschool GREENSBOROUGH',
  'li: Gurri Wanyarra Wellbeing Centre KANGAROO FLAT',
  'li: Huntly Primary School OSHC HUNTLY',
  'li: Kangaroo Ground Primary Combined OSHC KANGAROO GROUND',
  'li: Kent Road Uniting Church Preschool PASCOE VALE',
  'li: Mildura Waves Leisure Centre MILDURA',
  'li: New Futures Broadmeadows BROADMEADOWS',
  'li: New Futures Epping EPPING',
  'li: New Futures Occasional Child Care Centre COBURG',
  'li: Newbury Child and Community Centre CRAIGIEBURN',
  'li: Richmond West Afterschool Care and Vacation Care Program RICHMOND',
  'li: Romsey Occasional Care ROMSEY',
  'li: Roxburgh Park Community House & Children''s Centre ROXBURGH PARK',
  'li: St Andrews Uniting Church Kindergarten FAIRFIELD',
  'li: St Johns XXIII OSHClub Thomastown East THOMASTOWN',
  'li: Swan Hill Leisure Centre SWAN HILL',
  'li: Yarra Warra Preschool WARRANDYTE NORTH',
  'Schools closed',
  'li: The Department has',
  'not',
  'li: been advised of any school closures.',
  'TAFE',
  'li: The Department has notbeen advised of any TAFE closures.',
  'South-Western Victoria RegionEarly childhood services',
  'li: The Department has been advised of the following early childhood service closures:',
  'li: Aerotots Activity Centre WERRIBEE',
  'li: Altona Meadows Community Centre Occasional Care ALTONA MEADOWS',
  'li: Aquapulse Creche HOPPERS CROSSING',
  'li: Aspire Early Education Atherstone Estate STRATHTULLOH',
  'li: BASC Out of School Hours Care OCEAN GROVE',
  'li: Blackwood Street Neighbourhood House YARRAVILLE',
  'li: Bluewater Leisure Centre Creche COLAC',
  'li: Bridge Road Occasional Care MELTON SOUTH',
  'li: Camp Australia - Aireys Inlet Primary School OSHC AIREYS INLET',
  'li: Camp Australia - Baden Powell P-9 College Derrimut Heath Campus OSHC HOPPERS CROSSING',
  'li: Camp Australia - Christian College Geelong - Surf Coast Campus OSHC TORQUAY',
  'li: Camp Australia - Haileybury City Campus OSHC WEST MELBOURNE',
  'li: Camp Australia - Monmia Primary School OSHC KEILOR DOWNS',
  'li: Cana Catholic Primary OSHClub HILLSIDE',
  'li: Carranballac OSHClub Boardwalk POINT COOK',
  'li: Eagle Stadium WERRIBEE',
  'li: Energy Force Fitness Creche DRYSDALE',
  'li: Fernwood Fitness Sydenham SYDENHAM',
  'li: Fun 4 All Occasional Care Center WERRIBEE',
  'li: Genesis Maidstone MAIDSTONE',
  'li: Goodlife Essendon Child Minding ESSENDON',
  'li: Goodlife Geelong BELMONT',
  'li: Goodlife Point Cook POINT COOK',
  'li: Goodlife Taylors Lakes TAYLORS LAKES',
  'li: Hopetoun Child Care Service HOPETOUN',
  'li: Inspiring Minds Studio ESSENDON',
  'li: Iramoo Community Centre Occasional Child Care WYNDHAM VALE',
  'li: Keilor Basketball Stadium Creche KEILOR PARK',
  'li: Kensington Neighbourhood House Inc KENSINGTON',
  'li: Kids @ IAG Melbourne Vacation Care MELBOURNE',
  'li: Little Stars at Southern Cross DOCKLANDS',
  'li: Maribyrnong Aquatic Centre Occasional Child Care MARIBYRNONG',
  'li: Melton Central Occasional Care MELTON',
  'li: Melton Waves Leisure Centre Creche MELTON',
  'li: Myrniong Primary School Outside School Hours Care MYRNIONG',
  'li: New Futures Braybrook BRAYBROOK',
  'li: Pentland Afterschool Care Program DARLEY',
  'li: Shine Early Learning St Albans ST ALBANS',
  'li: Shuter Street Occasional Care MOONEE PONDS',
  'li: South Kingsville Community Centre SOUTH KINGSVILLE',
  'li: South Pacific Health Club Williamstown NEWPORT',
  'li: Springside Occasional Care CAROLINE SPRINGS',
  'li: Sunshine Leisure Centre SUNSHINE',
  'li: Warracknabeal Neighbourhood House Take a Break Child Care Centre WARRACKNABEAL',
  'li: Willaura Primary School OSHC ARARAT',
  'li: Winchelsea Primary School Winchelsea',
  'li: Yarraville Community Centre YARRAVILLE',
  'Schools closed',
  'li: The Department has been advised of the following school closures:',
  'Al Taqwa College, TRUGANINAParkville College (Malmsbury campus), PARKVILLETAFE',
  'li: The Department has not been advised of any TAFE closures.',
  'Bus service cancellations or alterations',
  'li: For Term 3 2020, schools bus services will continue to be provided to support student travel to schools where needed.',
  'li: With learning from home arrangements, local principals are authorised to cancel or suspend school buses where not required.',
  'li: Families and students are encouraged to contact their school directly to enquire if their bus service is continuing to be provided.',
  'li: Find more about bus services:',
  'School bus services in Term 3 (for schools)School bus services in Term 3 (for parents)VicRoadsRegional roadsBack to top',
  'li: Last Update: 08 September 2020',
  'In this section:- Example left hand nav using current codeHealth, wellbeing and safetyWebsite navigationFor parentsEarly childhoodSchoolsTAFE and trainingThe Department<li>                                 <a href="/about/educationstate">Education State                                                                  </a>                             </li>Support linksHelp in your languageAccessibilityPrivacyCopyrightDisclaimerContact linksContact usReport a website issue',
  'li: State Government of Victoria, Australia © 2019',
  'li: We respectfully acknowledge the Traditional Owners of country throughout Victoria and pay respect to the ongoing living cultures of First Peoples.',
  'li: Our website uses a free tool to translate into other languages. This tool is a guide and may not be accurate. For more, see:',
  'Information in your languageOld Search Code PlaceholdCludo Code for www.education.vic.gov.auHotjar Tracking Code for www.education.vic.gov.auSR-1181393'
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

Write-Host "Updated $($values.Count) rows in column A"